# Add a new, blank slide after the existing slide (becomes slide 2) —
# "Added Display blocks reference content".
$p = $ppt.ActivePresentation

# ppLayoutBlank = 12 -> corresponds to the deck's "Blank" slide layout
# (slideLayout7.xml), which is what a freshly inserted blank slide uses.
$newSlide = $p.Slides.Add(2, 12)
